$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 758.16455
$ws.Range("I15").Value = 758.16455
$ws.Range("K15").Value = 2274.49365
$ws.Range("M15").Value = -2105.49365

$ws.Range("H75").Value = 101799
$ws.Range("J75").Value = 101799
$ws.Range("L75").Value = 101799
$ws.Range("N75").Value = -103671

$ws.Range("H78").Value = 101799
$ws.Range("J78").Value = 101799
$ws.Range("L78").Value = 305397
$ws.Range("N78").Value = -314757

$ws.Range("H92").Value = 1889.2858
$ws.Range("I92").Value = 2276.8
$ws.Range("K92").Value = 2276.8
$ws.Range("M92").Value = -1028.8

$ws.Range("H96").Value = 2686.875
$ws.Range("I96").Value = 1946.7368
$ws.Range("K96").Value = 5840.2104
$ws.Range("M96").Value = -4467.2104

$ws.Range("H116").Value = 4869.0625
$ws.Range("I116").Value = 4878.3335
$ws.Range("J116").Value = 4857.143
$ws.Range("K116").Value = 4878.3335
$ws.Range("L116").Value = 4857.143
$ws.Range("M116").Value = -1436.3335
$ws.Range("N116").Value = -11741.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 28574000
$ws.Range("I61").Value = 37038220
$ws.Range("J61").Value = 7264.125
$ws.Range("K61").Value = 37038220
$ws.Range("L61").Value = 7264.125
$ws.Range("M61").Value = -37038008
$ws.Range("N61").Value = -7688.125

$ws.Range("H63").Value = 5264.926
$ws.Range("I63").Value = 3271.4736
$ws.Range("K63").Value = 3271.4736
$ws.Range("M63").Value = -2585.4736

$ws.Range("H66").Value = 5264.926
$ws.Range("I66").Value = 3271.4736
$ws.Range("K66").Value = 16357.368
$ws.Range("M66").Value = -12925.368

$ws.Range("H74").Value = 52692080
$ws.Range("I74").Value = 66742148
$ws.Range("J74").Value = 4326
$ws.Range("K74").Value = 66742148
$ws.Range("L74").Value = 4326
$ws.Range("M74").Value = -66741274
$ws.Range("N74").Value = -6074

$ws.Range("H77").Value = 52692080
$ws.Range("I77").Value = 66742148
$ws.Range("J77").Value = 4326
$ws.Range("K77").Value = 333710740
$ws.Range("L77").Value = 21630
$ws.Range("M77").Value = -333706372
$ws.Range("N77").Value = -30366

$ws.Range("H136").Value = 28574000
$ws.Range("I136").Value = 37038220
$ws.Range("J136").Value = 7264.125
$ws.Range("K136").Value = 111114660
$ws.Range("L136").Value = 21792.375
$ws.Range("M136").Value = -111112110
$ws.Range("N136").Value = -26892.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1124.5385
$ws.Range("I105").Value = 1100.8
$ws.Range("J105").Value = 1203.6666
$ws.Range("K105").Value = 1100.8
$ws.Range("L105").Value = 1203.6666
$ws.Range("M105").Value = 646.2
$ws.Range("N105").Value = -4697.6666

$ws.Range("H134").Value = 3087.276
$ws.Range("I134").Value = 3156.9614
$ws.Range("K134").Value = 9470.8842
$ws.Range("M134").Value = -6935.8842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16953314
$ws.Range("I31").Value = 3176.1428
$ws.Range("J31").Value = 41672264
$ws.Range("K31").Value = 3176.1428
$ws.Range("L31").Value = 41672264
$ws.Range("M31").Value = -2881.1428
$ws.Range("N31").Value = -41672854

$ws.Range("H34").Value = 16953314
$ws.Range("I34").Value = 3176.1428
$ws.Range("J34").Value = 41672264
$ws.Range("K34").Value = 3176.1428
$ws.Range("L34").Value = 41672264
$ws.Range("M34").Value = -2974.1428
$ws.Range("N34").Value = -41672668

$ws.Range("H132").Value = 4137
$ws.Range("I132").Value = 3987.4707
$ws.Range("J132").Value = 4772.5
$ws.Range("K132").Value = 11962.4121
$ws.Range("L132").Value = 14317.5
$ws.Range("M132").Value = -9432.4121
$ws.Range("N132").Value = -19377.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 785.1818
$ws.Range("I2").Value = 1182
$ws.Range("K2").Value = 7092
$ws.Range("M2").Value = -6979

$ws.Range("H92").Value = 499
$ws.Range("I92").Value = 499
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1497
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -249
$ws.Range("N92").ClearContents()

$ws.Range("H93").Value = 549.5
$ws.Range("I93").Value = 100
$ws.Range("J93").Value = 999
$ws.Range("K93").Value = 300
$ws.Range("L93").Value = 2997
$ws.Range("M93").Value = 1572
$ws.Range("N93").Value = -6741

$ws.Range("H96").Value = 2025
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H121").Value = 1134
$ws.Range("J121").Value = 1999.75
$ws.Range("L121").Value = 5999.25
$ws.Range("N121").Value = -8619.25

$ws.Range("H131").Value = 32892.11
$ws.Range("J131").Value = 5362.76
$ws.Range("L131").Value = 16088.28
$ws.Range("N131").Value = -26168.28

$ws.Range("H132").Value = 1484172.5
$ws.Range("I132").Value = 2361.9048
$ws.Range("J132").Value = 2780756.8
$ws.Range("K132").Value = 21257.1432
$ws.Range("L132").Value = 25026811.2
$ws.Range("M132").Value = -18727.1432
$ws.Range("N132").Value = -25031871.2

$ws.Range("H133").Value = 11160.556
$ws.Range("J133").Value = 19504
$ws.Range("L133").Value = 58512
$ws.Range("N133").Value = -68632

$ws.Range("H134").Value = 3058.2666
$ws.Range("J134").Value = 19499
$ws.Range("L134").Value = 58497
$ws.Range("N134").Value = -68637

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3200.4614
$ws.Range("I80").Value = 2689.625
$ws.Range("K80").Value = 2689.625
$ws.Range("M80").Value = -1691.625

$ws.Range("H83").Value = 3200.4614
$ws.Range("I83").Value = 2689.625
$ws.Range("K83").Value = 13448.125
$ws.Range("M83").Value = -8456.125

$ws.Range("H122").Value = 1828.8422
$ws.Range("I122").Value = 1582.4615
$ws.Range("J122").Value = 2362.6667
$ws.Range("K122").Value = 4747.3845
$ws.Range("L122").Value = 7088.000100000001
$ws.Range("M122").Value = -2297.3845
$ws.Range("N122").Value = -11988.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 542
$ws.Range("I55").Value = 300.6
$ws.Range("J55").Value = 904.1
$ws.Range("K55").Value = 300.6
$ws.Range("L55").Value = 904.1
$ws.Range("M55").Value = -127.6
$ws.Range("N55").Value = -1250.1

$ws.Range("H100").Value = 3528.2222
$ws.Range("I100").Value = 2225
$ws.Range("K100").Value = 2225
$ws.Range("M100").Value = -1684

$ws.Range("H122").Value = 3318.074
$ws.Range("I122").Value = 2885.3809
$ws.Range("J122").Value = 4832.5
$ws.Range("K122").Value = 8656.1427
$ws.Range("L122").Value = 14497.5
$ws.Range("M122").Value = -6206.1427
$ws.Range("N122").Value = -19397.5

$ws.Range("H132").Value = 83334300
$ws.Range("I132").Value = 1077.7368
$ws.Range("J132").Value = 400000540
$ws.Range("K132").Value = 3233.2104
$ws.Range("L132").Value = 1200001620
$ws.Range("M132").Value = -703.2103999999999
$ws.Range("N132").Value = -1200006680

$ws.Range("H136").Value = 1910.5962
$ws.Range("I136").Value = 1910.5962
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5731.7886
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3181.7886
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7393.9443
$ws.Range("I62").Value = 6940.8
$ws.Range("J62").Value = 7568.231
$ws.Range("K62").Value = 6940.8
$ws.Range("L62").Value = 7568.231
$ws.Range("M62").Value = -6316.8
$ws.Range("N62").Value = -8816.231

$ws.Range("H65").Value = 7393.9443
$ws.Range("I65").Value = 6940.8
$ws.Range("J65").Value = 7568.231
$ws.Range("K65").Value = 34704
$ws.Range("L65").Value = 37841.155
$ws.Range("M65").Value = -31584
$ws.Range("N65").Value = -44081.155

$ws.Range("H81").Value = 1244.25
$ws.Range("I81").Value = 887.5
$ws.Range("K81").Value = 1775
$ws.Range("M81").Value = -714

$ws.Range("H84").Value = 1244.25
$ws.Range("I84").Value = 887.5
$ws.Range("K84").Value = 8875
$ws.Range("M84").Value = -3571

$ws.Range("H132").Value = 4209.3335
$ws.Range("I132").Value = 4080.0698
$ws.Range("K132").Value = 12240.2094
$ws.Range("M132").Value = -9710.2094
